# Auto-generated edit script
# Applies value updates to the Typhon_Profits leve-profit tables across sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H3").Value = 29800
$ws.Range("J3").Value = 29800
$ws.Range("L3").Value = 29800
$ws.Range("N3").Value = -30028
$ws.Range("H28").Value = 181.75
$ws.Range("I28").Value = 189.18182
$ws.Range("J28").Value = 100
$ws.Range("K28").Value = 189.18182
$ws.Range("L28").Value = 100
$ws.Range("M28").Value = 295.81818
$ws.Range("N28").Value = -1070
$ws.Range("H98").Value = 369.88235
$ws.Range("I98").Value = 386.75
$ws.Range("J98").Value = 100
$ws.Range("K98").Value = 386.75
$ws.Range("L98").Value = 100
$ws.Range("M98").Value = 1111.25
$ws.Range("N98").Value = -3096
$ws.Range("H100").Value = 71432300
$ws.Range("I100").Value = 250001740
$ws.Range("K100").Value = 250001740
$ws.Range("M100").Value = -250001199
$ws.Range("H102").Value = 29800
$ws.Range("J102").Value = 29800
$ws.Range("L102").Value = 29800
$ws.Range("N102").Value = -36290
$ws.Range("H122").Value = 369.88235
$ws.Range("I122").Value = 386.75
$ws.Range("J122").Value = 100
$ws.Range("K122").Value = 1160.25
$ws.Range("L122").Value = 300
$ws.Range("M122").Value = 1289.75
$ws.Range("N122").Value = -5200
$ws.Range("H129").Value = 173240.19
$ws.Range("J129").Value = 182669.66
$ws.Range("L129").Value = 548008.98
$ws.Range("N129").Value = -558008.98
$ws.Range("H135").Value = 16135064
$ws.Range("J135").Value = 45470000
$ws.Range("L135").Value = 409230000
$ws.Range("N135").Value = -409235070

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H110").Value = 684.17645
$ws.Range("I110").Value = 617.25
$ws.Range("J110").Value = 844.8
$ws.Range("K110").Value = 617.25
$ws.Range("L110").Value = 844.8
$ws.Range("M110").Value = 1427.75
$ws.Range("N110").Value = -4934.8

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 2214
$ws.Range("I75").Value = 2214
$ws.Range("K75").Value = 2214
$ws.Range("M75").Value = -1278
$ws.Range("H78").Value = 2214
$ws.Range("I78").Value = 2214
$ws.Range("K78").Value = 6642
$ws.Range("M78").Value = -1962
$ws.Range("H86").Value = 1734.7106
$ws.Range("I86").Value = 1441.5518
$ws.Range("J86").Value = 2679.3333
$ws.Range("K86").Value = 1441.5518
$ws.Range("L86").Value = 2679.3333
$ws.Range("M86").Value = -318.5518
$ws.Range("N86").Value = -4925.3333
$ws.Range("H89").Value = 1734.7106
$ws.Range("I89").Value = 1441.5518
$ws.Range("J89").Value = 2679.3333
$ws.Range("K89").Value = 7207.759
$ws.Range("L89").Value = 13396.6665
$ws.Range("M89").Value = -1591.759
$ws.Range("N89").Value = -24628.6665
$ws.Range("H94").Value = 1704.5172
$ws.Range("I94").Value = 1509.591
$ws.Range("J94").Value = 2317.1428
$ws.Range("K94").Value = 1509.591
$ws.Range("L94").Value = 2317.1428
$ws.Range("M94").Value = -1058.591
$ws.Range("N94").Value = -3219.1428

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 13568.519
$ws.Range("I31").Value = 18791.176
$ws.Range("J31").Value = 4690
$ws.Range("K31").Value = 18791.176
$ws.Range("L31").Value = 4690
$ws.Range("M31").Value = -18496.176
$ws.Range("N31").Value = -5280
$ws.Range("H34").Value = 13568.519
$ws.Range("I34").Value = 18791.176
$ws.Range("J34").Value = 4690
$ws.Range("K34").Value = 18791.176
$ws.Range("L34").Value = 4690
$ws.Range("M34").Value = -18589.176
$ws.Range("N34").Value = -5094
$ws.Range("H58").Value = 15348.286
$ws.Range("I58").Value = 1190.4286
$ws.Range("J58").Value = 36585.07
$ws.Range("K58").Value = 1190.4286
$ws.Range("L58").Value = 36585.07
$ws.Range("M58").Value = -987.4286
$ws.Range("N58").Value = -36991.07
$ws.Range("H110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("L110").Value = ""
$ws.Range("N110").Value = 0
$ws.Range("H122").Value = 1762.5714
$ws.Range("I122").Value = 2272
$ws.Range("J122").Value = 1083.3334
$ws.Range("K122").Value = 6816
$ws.Range("L122").Value = 3250.0002
$ws.Range("M122").Value = -4366
$ws.Range("N122").Value = -8150.0002
$ws.Range("H136").Value = 15348.286
$ws.Range("I136").Value = 1190.4286
$ws.Range("J136").Value = 36585.07
$ws.Range("K136").Value = 3571.2858
$ws.Range("L136").Value = 109755.21
$ws.Range("M136").Value = -1021.2858
$ws.Range("N136").Value = -114855.21

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 726.38
$ws.Range("J131").Value = 749.337
$ws.Range("L131").Value = 2248.011
$ws.Range("N131").Value = -12328.011
$ws.Range("H136").Value = 2444.8
$ws.Range("I136").Value = 1806.25
$ws.Range("K136").Value = 5418.75
$ws.Range("M136").Value = -318.75
$ws.Range("H140").Value = 2477.5
$ws.Range("I140").Value = 2205
$ws.Range("K140").Value = 6615
$ws.Range("M140").Value = -1435

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H52").Value = 25005250
$ws.Range("J52").Value = 25005250
$ws.Range("L52").Value = 25005250
$ws.Range("N52").Value = -25005768
$ws.Range("H70").Value = 13859.3
$ws.Range("I70").Value = 19899.5
$ws.Range("J70").Value = 4799
$ws.Range("K70").Value = 19899.5
$ws.Range("L70").Value = 4799
$ws.Range("M70").Value = -19629.5
$ws.Range("N70").Value = -5339
$ws.Range("H73").Value = 13859.3
$ws.Range("I73").Value = 19899.5
$ws.Range("J73").Value = 4799
$ws.Range("K73").Value = 19899.5
$ws.Range("L73").Value = 4799
$ws.Range("M73").Value = -18963.5
$ws.Range("N73").Value = -6671
$ws.Range("H80").Value = 3738.5908
$ws.Range("J80").Value = 3927.8667
$ws.Range("L80").Value = 3927.8667
$ws.Range("N80").Value = -5923.8667
$ws.Range("H83").Value = 3738.5908
$ws.Range("J83").Value = 3927.8667
$ws.Range("L83").Value = 19639.3335
$ws.Range("N83").Value = -29623.3335
$ws.Range("H102").Value = 21741502
$ws.Range("I102").Value = 29414336
$ws.Range("K102").Value = 29414336
$ws.Range("M102").Value = -29412714

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 1672.1904
$ws.Range("I93").Value = 1539.3334
$ws.Range("K93").Value = 1539.3334
$ws.Range("M93").Value = -291.3334
$ws.Range("H132").Value = 2145.5667
$ws.Range("I132").Value = 1539.5
$ws.Range("K132").Value = 4618.5
$ws.Range("M132").Value = -2088.5
